# Added implementation of MSM measure.
# The "methodNumberOfLines" sheet listed a row for every method, including
# no-arg/implicit constructors. Those constructor rows are removed here,
# shifting the remaining method rows up (Excel's Rows.Delete shift-up
# semantics), which also drops the now-unused "4" shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("methodNumberOfLines")

# Rows (1-based, as currently laid out) holding constructor entries:
#   row 2  -> OrderManageService(CustomerRepository, KafkaTemplate)  [4 lines]
#   row 5  -> PaymentAppTest()                                       [1 line]
#   row 7  -> PaymentComponentTests()                                [1 line]
#   row 11 -> PaymentApp()                                           [1 line]
#   row 24 -> KafkaContainerDevMode()                                [1 line]
# Delete from the bottom up so earlier row numbers stay valid.
$rowsToDelete = @(24, 11, 7, 5, 2)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
